$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 2; this shifts the existing "Normal" (row 2)
# and "Model" (row 3) rows down to rows 3 and 4.
$ws.Rows("2:2").Insert()

# The insert leaves the shifted rows carrying a row-wide style/customFormat
# attribute. Strip that blanket row formatting while keeping each row's
# per-cell style intact, by round-tripping the formats through a scratch
# range far away from the used data.
$ws.Range("A3:K3").Copy()
$ws.Range("A20:K20").PasteSpecial(-4122)
$ws.Rows("3:3").ClearFormats()
$ws.Range("A20:K20").Copy()
$ws.Range("A3:K3").PasteSpecial(-4122)
$ws.Range("A20:K20").Clear()

$ws.Range("A4:K4").Copy()
$ws.Range("A20:K20").PasteSpecial(-4122)
$ws.Rows("4:4").ClearFormats()
$ws.Range("A20:K20").Copy()
$ws.Range("A4:K4").PasteSpecial(-4122)
$ws.Range("A20:K20").Clear()

# Give the new row 2 the same (green) fill style as the "Model" row, then
# fill in the "LO2" label and its values.
$ws.Range("A4:F4").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$ws.Range("A2").Value = "LO2"
$ws.Range("B2").Value = 2.01
$ws.Range("C2").Value = 1.64
$ws.Range("D2").Value = 1.37
$ws.Range("E2").Value = 1.25
$ws.Range("F2").Value = 0.83

# Row 1 (header) gets a taller, explicit height.
$ws.Rows("1:1").RowHeight = 18.75

# Restore the selection to match the author's saved view state.
$ws.Range("G12").Select()
